# The "mango-mochi-project" BLS workbook: drop the 2000-2002 rows because
# there is no Asian unemployment data before 2003, shifting the whole
# 2003-2020 series up by three rows.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("BLS Data Series")
$metaSheet = $wb.Worksheets.Item("Sheet1")

$dataSheet.Activate()

# --- 1. Remove the obsolete 2000, 2001 and 2002 rows -----------------------
# Row 2 = 2000, Row 3 = 2001, Row 4 = 2002; deleting them shifts 2003..2020
# up into rows 2..19.
$dataSheet.Rows("2:4").Delete()

# --- 2. Update the frozen-pane selection on the data sheet ------------------
$dataSheet.Range("A2:XFD4").Select() | Out-Null

# --- 3. Grow the two wrapped-text description rows on the metadata sheet ---
# ("Labor force status:" and "Type of data:" explanations now need extra
# height once re-wrapped.)
$metaSheet.Rows("7:7").RowHeight = 28
$metaSheet.Rows("8:8").RowHeight = 28
